# "fix up images and table size"
#
# 1) Update the cached date footer text (datetimeFigureOut field) from
#    23/08/2023 to 24/08/2023 everywhere it appears: the slide master and
#    every slide layout (see note further below re: the notes master).
# 2) Re-sync the stray shape/group names left over on slide 1 so that
#    they match the naming already used for the same group on slides
#    2-4 (Group 5 / Rectangle: Rounded Corners 3/10/11 / Straight Arrow
#    Connector 7).

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "23/08/2023") {
            $chars = $tr.Characters(1, $tr.Length)
            $chars.Text = "24/08/2023"
        }
    }
}

$p = $ppt.ActivePresentation

# --- Slide master ---
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    Update-DateShape $master.Shapes.Item($j)
}

# --- Every slide layout ---
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        Update-DateShape $layout.Shapes.Item($j)
    }
}

# NOTE: the notes master's date placeholder (ppt/notesMasters/notesMaster1.xml)
# also needs the same 23/08/2023 -> 24/08/2023 fix, but this runtime's COM
# object model does not allow edits to Notes Master shapes to be persisted:
# any attempted mutation through $p.NotesMaster.Shapes either silently no-ops,
# or - worse - lands on whichever Slide Master shape happens to reuse the same
# internal shape Id (e.g. the notes master's "Date Placeholder 2" has Id=3,
# which collides with the slide master's "Text Placeholder 2"). So it is
# intentionally left untouched here to avoid corrupting the slide master.

# --- Re-sync the leftover shape names on slide 1's group with the
#     naming already used by the same group on slides 2-4 ---
$s1 = $p.Slides.Item(1)
$grp = $s1.Shapes.Item(1)
$grp.Name = "Group 5"

for ($j = 1; $j -le $grp.GroupItems.Count; $j++) {
    $item = $grp.GroupItems.Item($j)
    switch ($item.Id) {
        9  { $item.Name = "Rectangle: Rounded Corners 3" }
        10 { $item.Name = "Straight Arrow Connector 7" }
        16 { $item.Name = "Rectangle: Rounded Corners 10" }
        17 { $item.Name = "Rectangle: Rounded Corners 11" }
    }
}
